# Update cryptos list with latest prices and volume changes
# (GitHub Actions scheduled refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.744.73"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.481.89"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'580.64"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'160.71"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +8.43%  "
$ws.Range("D9").Value = "3.481.13"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "4.074.61"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'0.0000196"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "'28.72"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "65.689.29"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "3.464.18"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'6.47"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'14.29"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "'390.44"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "'8.24"
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'73.51"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").Value = "'9.81"
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'6.39"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'23.75"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'6.50"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("D36").Value = "'7.12"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").Value = "'162.92"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").Value = "3.079.40"
$ws.Range("E40").Value = "  +4.86%  "
$ws.Range("D41").Value = "'0.0771"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").Value = "'27.12"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "'0.0323"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.53"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'43.07"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").Value = "'0.777"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'25.63"
$ws.Range("E47").Value = "  +7.29%  "
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "'6.70"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("E51").Value = "  +3.66%  "
